# Generate Report for Handback
# Updates row 7 (bc729f1e-4578-4c1a-8788-a28f9ff4835f) on both the zh-cn and
# de-de sheets: it now has a Latest Target File / Latest Handback File /
# Latest Handback DateTime, plus an Error Detail explaining the handback
# file version is stale, and a new hyperlink on the Latest Handback File
# cell (column I).

$wb = $excel.ActiveWorkbook

$staleMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0d8ba40a0085b381e5297cc98a43e8db488aaf4d/e2e/bc729f1e-4578-4c1a-8788-a28f9ff4835f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f358937fe1ea90292976d004d9ba57866fa0d7b/e2e/bc729f1e-4578-4c1a-8788-a28f9ff4835f.md."
$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f358937fe1ea90292976d004d9ba57866fa0d7b/e2e/bc729f1e-4578-4c1a-8788-a28f9ff4835f.md"
$handbackMd = "bc729f1e-4578-4c1a-8788-a28f9ff4835f.md"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("I7").Value = $handbackMd
$wsZh.Range("J7").Value = "bc729f1e-4578-4c1a-8788-a28f9ff4835f.c289895ebddd3147427cf3e85fc8bb8a8203816e.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-09-04 17:01:05"
$wsZh.Range("P7").Value = $staleMessage
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $handbackUrl, "", "", $handbackMd)

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("I7").Value = $handbackMd
$wsDe.Range("J7").Value = "bc729f1e-4578-4c1a-8788-a28f9ff4835f.c289895ebddd3147427cf3e85fc8bb8a8203816e.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-04 17:01:15"
$wsDe.Range("P7").Value = $staleMessage
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $handbackUrl, "", "", $handbackMd)
